$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 229, shifting existing rows 229:334 down to 230:335
$ws.Rows("229:229").Insert()

# Populate the newly inserted row 229 with the new record
$ws.Range("A229").Value = 11
$ws.Range("B229").Value = "Vega Monumental Concepción"
$ws.Range("C229").Value = "Bíobío"
$ws.Range("D229").Value = 44875
$ws.Range("E229").Value = 8
$ws.Range("F229").Value = 100112008
$ws.Range("G229").Value = "Coliflor"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 2200
$ws.Range("K229").Value = 700
$ws.Range("L229").Value = 800
$ws.Range("M229").Value = 745
$ws.Range("N229").Value = "$/unidad"
$ws.Range("O229").Value = "Región Metropolitana"
$ws.Range("P229").Value = 745
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"
